# Update the VIC second-doses sheet with data through 8 November.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert five new rows at the top (newest date first) ----------------
$newRows = @(
    @{ Date = 44508; Value = 4774679 },
    @{ Date = 44507; Value = 4761123 },
    @{ Date = 44506; Value = 4733663 },
    @{ Date = 44505; Value = 4689787 },
    @{ Date = 44504; Value = 4643853 }
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $ws.Rows.Item(2).Insert()
    $r = 2
    $ws.Cells.Item($r, 1).Value = $newRows[$i].Date
    $ws.Cells.Item($r, 2).Value = $newRows[$i].Value
    $ws.Cells.Item($r, 1).Style = $ws.Cells.Item($r + 1, 1).Style
    $ws.Cells.Item($r, 2).Style = $ws.Cells.Item($r + 1, 2).Style
    $ws.Rows.Item($r).RowHeight = 18
}

# --- 2. Drop the old trailing (empty, style-only) row -----------------------
# After the 5-row insert above, the old empty row (228) now lives at 233.
$ws.Rows.Item(233).Delete()

# --- 3. Re-anchor the "Donate" picture up by one row ------------------------
$shp = $ws.Shapes.Item(1)
$shp.Top = 4066

# --- 4. Selection ------------------------------------------------------------
$ws.Range("E10").Select()
